$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 25646214
$ws.Range("I132").Value = 40006244
$ws.Range("J132").Value = 3308.3572
$ws.Range("K132").Value = 120018732
$ws.Range("L132").Value = 9925.071599999999
$ws.Range("M132").Value = -120016202
$ws.Range("N132").Value = -14985.0716
# Row 140
$ws.Range("H140").Value = 48611.332
$ws.Range("J140").Value = 49262.145
$ws.Range("L140").Value = 49262.145
$ws.Range("N140").Value = -59622.145

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1805.7142
$ws.Range("I45").Value = 1735
$ws.Range("K45").Value = 1735
$ws.Range("M45").Value = -1358
# Row 61
$ws.Range("H61").Value = 994.02856
$ws.Range("I61").Value = 840.95654
$ws.Range("J61").Value = 1287.4166
$ws.Range("K61").Value = 840.95654
$ws.Range("L61").Value = 1287.4166
$ws.Range("M61").Value = -628.95654
$ws.Range("N61").Value = -1711.4166
# Row 80
$ws.Range("H80").Value = 31516.666
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 31516.666
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 31516.666
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -33512.666
# Row 83
$ws.Range("H83").Value = 31516.666
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 31516.666
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 94549.99800000001
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -104533.998
# Row 136
$ws.Range("H136").Value = 994.02856
$ws.Range("I136").Value = 840.95654
$ws.Range("J136").Value = 1287.4166
$ws.Range("K136").Value = 2522.86962
$ws.Range("L136").Value = 3862.2498
$ws.Range("M136").Value = 27.13038000000006
$ws.Range("N136").Value = -8962.2498

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1816
$ws.Range("I86").Value = 1768.6666
$ws.Range("J86").Value = 1863.3334
$ws.Range("K86").Value = 1768.6666
$ws.Range("L86").Value = 1863.3334
$ws.Range("M86").Value = -645.6666
$ws.Range("N86").Value = -4109.3334
# Row 89
$ws.Range("H89").Value = 1816
$ws.Range("I89").Value = 1768.6666
$ws.Range("J89").Value = 1863.3334
$ws.Range("K89").Value = 8843.333000000001
$ws.Range("L89").Value = 9316.666999999999
$ws.Range("M89").Value = -3227.333000000001
$ws.Range("N89").Value = -20548.667
# Row 107
$ws.Range("H107").Value = 2138.6667
$ws.Range("I107").Value = 2084.4
$ws.Range("J107").Value = 2206.5
$ws.Range("K107").Value = 2084.4
$ws.Range("L107").Value = 2206.5
$ws.Range("M107").Value = -164.4000000000001
$ws.Range("N107").Value = -6046.5
# Row 134
$ws.Range("H134").Value = 2099.86
$ws.Range("I134").Value = 1032.6052
$ws.Range("J134").Value = 5479.5
$ws.Range("K134").Value = 3097.8156
$ws.Range("L134").Value = 16438.5
$ws.Range("M134").Value = -562.8155999999999
$ws.Range("N134").Value = -21508.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 10206526
$ws.Range("I31").Value = 1319.7812
$ws.Range("J31").Value = 29416324
$ws.Range("K31").Value = 1319.7812
$ws.Range("L31").Value = 29416324
$ws.Range("M31").Value = -1024.7812
$ws.Range("N31").Value = -29416914
# Row 34
$ws.Range("H34").Value = 10206526
$ws.Range("I34").Value = 1319.7812
$ws.Range("J34").Value = 29416324
$ws.Range("K34").Value = 1319.7812
$ws.Range("L34").Value = 29416324
$ws.Range("M34").Value = -1117.7812
$ws.Range("N34").Value = -29416728
# Row 48
$ws.Range("H48").Value = 46683.668
$ws.Range("J48").Value = 46683.668
$ws.Range("L48").Value = 46683.668
$ws.Range("N48").Value = -47635.668
# Row 87
$ws.Range("H87").Value = 22250
$ws.Range("J87").Value = 22250
$ws.Range("L87").Value = 22250
$ws.Range("N87").Value = -24622
# Row 90
$ws.Range("H90").Value = 22250
$ws.Range("J90").Value = 22250
$ws.Range("L90").Value = 66750
$ws.Range("N90").Value = -78606
# Row 94
$ws.Range("H94").Value = 1258.4
$ws.Range("I94").Value = 738.7143
$ws.Range("J94").Value = 1538.2307
$ws.Range("K94").Value = 738.7143
$ws.Range("L94").Value = 1538.2307
$ws.Range("M94").Value = -287.7143
$ws.Range("N94").Value = -2440.2307
# Row 122
$ws.Range("H122").Value = 2692.0908
$ws.Range("I122").Value = 1302.6
$ws.Range("K122").Value = 3907.8
$ws.Range("M122").Value = -1457.8
# Row 132
$ws.Range("H132").Value = 3066.9167
$ws.Range("I132").Value = 2564.8386
$ws.Range("K132").Value = 7694.5158
$ws.Range("M132").Value = -5164.5158
# Row 134
$ws.Range("H134").Value = 4591.1353
$ws.Range("I134").Value = 5751.7617
$ws.Range("J134").Value = 3067.8125
$ws.Range("K134").Value = 17255.2851
$ws.Range("L134").Value = 9203.4375
$ws.Range("M134").Value = -14720.2851
$ws.Range("N134").Value = -14273.4375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 1285.5714
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 1399.8
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 4199.4
$ws.Range("M86").Value = -1814
$ws.Range("N86").Value = -6571.4
# Row 89
$ws.Range("H89").Value = 1285.5714
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 1399.8
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 12598.2
$ws.Range("M89").Value = -3072
$ws.Range("N89").Value = -24454.2
# Row 98
$ws.Range("H98").Value = 259.4375
$ws.Range("I98").Value = 93
$ws.Range("J98").Value = 314.91666
$ws.Range("K98").Value = 279
$ws.Range("L98").Value = 944.7499799999999
$ws.Range("M98").Value = 1219
$ws.Range("N98").Value = -3940.74998
# Row 113
$ws.Range("H113").Value = 857.4815
$ws.Range("I113").Value = 709.6111
$ws.Range("J113").Value = 1153.2222
$ws.Range("K113").Value = 2128.8333
$ws.Range("L113").Value = 3459.6666
$ws.Range("M113").Value = 41.16670000000022
$ws.Range("N113").Value = -7799.6666

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 17285.428
$ws.Range("J39").Value = 17285.428
$ws.Range("L39").Value = 17285.428
$ws.Range("N39").Value = -18349.428
# Row 96
$ws.Range("H96").Value = 14990
$ws.Range("J96").Value = 14990
$ws.Range("L96").Value = 14990
$ws.Range("N96").Value = -20482

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 41
$ws.Range("H41").Value = 10030
$ws.Range("I41").Value = 10030
$ws.Range("K41").Value = 10030
$ws.Range("M41").Value = -9592
# Row 43
$ws.Range("H43").Value = 29332.334
$ws.Range("J43").Value = 29332.334
$ws.Range("L43").Value = 29332.334
$ws.Range("N43").Value = -29718.334
# Row 50
$ws.Range("H50").Value = 34941
$ws.Range("J50").Value = 34941
$ws.Range("L50").Value = 34941
$ws.Range("N50").Value = -36215
# Row 54
$ws.Range("H54").Value = 35027
$ws.Range("J54").Value = 35027
$ws.Range("L54").Value = 35027
$ws.Range("N54").Value = -36315
# Row 75
$ws.Range("H75").Value = 38000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
# Row 78
$ws.Range("H78").Value = 38000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
# Row 100
$ws.Range("H100").Value = 2100.077
$ws.Range("I100").Value = 1922.5555
$ws.Range("J100").Value = 2499.5
$ws.Range("K100").Value = 1922.5555
$ws.Range("L100").Value = 2499.5
$ws.Range("M100").Value = -1381.5555
$ws.Range("N100").Value = -3581.5
# Row 123
$ws.Range("H123").Value = 29993.084
$ws.Range("J123").Value = 29993.084
$ws.Range("L123").Value = 29993.084
$ws.Range("N123").Value = -39793.084
# Row 132
$ws.Range("H132").Value = 8615.166999999999
$ws.Range("I132").Value = 8611.333000000001
$ws.Range("J132").Value = 8623.6
$ws.Range("K132").Value = 25833.999
$ws.Range("L132").Value = 25870.8
$ws.Range("M132").Value = -23303.999
$ws.Range("N132").Value = -30930.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 1166.6666
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1166.6666
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1166.6666
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -1502.6666
# Row 95
$ws.Range("H95").Value = 33562.668
$ws.Range("J95").Value = 33562.668
$ws.Range("L95").Value = 33562.668
$ws.Range("N95").Value = -39054.668
# Row 132
$ws.Range("H132").Value = 5849460.5
$ws.Range("I132").Value = 1073.762
$ws.Range("J132").Value = 22224942
$ws.Range("K132").Value = 3221.286
$ws.Range("L132").Value = 66674826
$ws.Range("M132").Value = -691.2860000000001
$ws.Range("N132").Value = -66679886
# Row 136
$ws.Range("H136").Value = 1091.5111
$ws.Range("I136").Value = 562.3939
$ws.Range("J136").Value = 2546.5833
$ws.Range("K136").Value = 1687.1817
$ws.Range("L136").Value = 7639.749899999999
$ws.Range("M136").Value = 862.8182999999999
$ws.Range("N136").Value = -12739.7499
